$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 7) describing a new logged test-mail interaction.
$ws.Range("A7").Value = "Testmail #6: Hebben we EcoPro-700 nog op voorraad?"
$ws.Range("B7").Value = "Beste klant,`nBedankt voor uw vraag. Op dit moment hebben we EcoPro-700 niet op voorraad. We verwachten binnenkort nieuwe voorraad binnen te krijgen. Mocht u nog vragen hebben of een pre-order willen plaatsen, neem dan gerust contact met ons op.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$ws.Range("C7").Value = "Hebben we EcoPro-700 nog op voorraad?"
$ws.Range("D7").Value = "mailmind.test@zohomail.eu"
$ws.Range("E7").Value = "Productinformatie"
$ws.Range("F7").Value = "2025-07-31 21:34:42"
$ws.Range("G7").Value = "Ja"
$ws.Range("H7").Value = "Nee"
$ws.Range("I7").Value = "Ja"
$ws.Range("J7").Value = "Nee"

# Keep the row height at its default (content-measured) state, matching the
# other data rows, instead of leaving an explicit autofit-computed height.
$ws.Rows.Item(7).AutoFit()
